$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New final state for rows 8-17 (columns A-E):
# row, A(idx), B(name), C(from_bus), D(to_bus), E(in_service)
$rows = @(
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16,  9, $true),
    @(10, 8,  "extr1",  5, 12, $true),
    @(11, 9,  "extr2",  5,  9, $true),
    @(12, 10, "extr3", 10, 11, $true),
    @(13, 11, "extr4",  7,  8, $true),
    @(14, 12, "extr5",  9, 11, $false),
    @(15, 13, "extr6",  7, 11, $true),
    @(16, 14, "extr7",  5,  7, $false),
    @(17, 15, "extr8",  8,  5, $true)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value2 = $r[1]
    $ws.Cells.Item($rowNum, 2).Value2 = $r[2]
    $ws.Cells.Item($rowNum, 3).Value2 = $r[3]
    $ws.Cells.Item($rowNum, 4).Value2 = $r[4]
    $ws.Cells.Item($rowNum, 5).Value2 = $r[5]
}

# Keep the bold/centered/bordered style used by the other id cells in column A
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
